$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 5
$ws.Range("F6").Value = 2413
$ws.Range("F8").Value = 1783
$ws.Range("F9").Value = 3028
$ws.Range("F10").Value = 182
$ws.Range("F11").Value = 4511
$ws.Range("F12").Value = 399
$ws.Range("F13").Value = 223
$ws.Range("F14").Value = 140
$ws.Range("F15").Value = 569
$ws.Range("F16").Value = 268
$ws.Range("F18").Value = 243
$ws.Range("F20").Value = 114
$ws.Range("F21").Value = 314
$ws.Range("F22").Value = 4549
$ws.Range("F24").Value = 4073
$ws.Range("F25").Value = 1147
$ws.Range("F26").Value = 220
$ws.Range("F27").Value = 593
$ws.Range("F28").Value = 4387
$ws.Range("F29").Value = 95
$ws.Range("F30").Value = 649
$ws.Range("F31").Value = 617
$ws.Range("F32").Value = 595
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 34
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1046
$ws.Range("F4").Value = 23
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1046
$ws.Range("F5").Value = 23
$ws.Range("F6").Value = 5
$ws.Range("F9").Value = 2413
$ws.Range("F11").Value = 1783
$ws.Range("F13").Value = 3028
$ws.Range("F14").Value = 182
$ws.Range("F15").Value = 4511
$ws.Range("F16").Value = 399
$ws.Range("F17").Value = 223
$ws.Range("F18").Value = 140
$ws.Range("F19").Value = 569
$ws.Range("F20").Value = 268
$ws.Range("F22").Value = 243
$ws.Range("F25").Value = 114
$ws.Range("F26").Value = 314
$ws.Range("F27").Value = 4549
$ws.Range("F29").Value = 4073
$ws.Range("F30").Value = 1147
$ws.Range("F31").Value = 220
$ws.Range("F32").Value = 593
$ws.Range("F33").Value = 4387
$ws.Range("F34").Value = 95
$ws.Range("F35").Value = 649
$ws.Range("F36").Value = 617
$ws.Range("F37").Value = 595
$ws.Range("F39").Value = 34
